$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 40
$ws.Range("H40").Value = 6449.727
$ws.Range("I40").Value = 10669
$ws.Range("J40").Value = 2933.6667
$ws.Range("K40").Value = 10669
$ws.Range("L40").Value = 2933.6667
$ws.Range("M40").Value = -10494
$ws.Range("N40").Value = -3283.6667
# Row 129
$ws.Range("H129").Value = 834.64
$ws.Range("I129").Value = 501.77777
$ws.Range("J129").Value = 1021.875
$ws.Range("K129").Value = 1505.33331
$ws.Range("L129").Value = 3065.625
$ws.Range("M129").Value = 3494.66669
$ws.Range("N129").Value = -13065.625
# Row 132
$ws.Range("H132").Value = 145021.45
$ws.Range("I132").Value = 2092.2615
$ws.Range("J132").Value = 2003101
$ws.Range("K132").Value = 6276.7845
$ws.Range("L132").Value = 6009303
$ws.Range("M132").Value = -3746.7845
$ws.Range("N132").Value = -6014363
# Row 137
$ws.Range("H137").Value = 6251.5
$ws.Range("I137").Value = 426
$ws.Range("J137").Value = 6898.778
$ws.Range("K137").Value = 1278
$ws.Range("L137").Value = 20696.334
$ws.Range("M137").Value = 1272
$ws.Range("N137").Value = -25796.334
# Row 138
$ws.Range("H138").Value = 2077.74
$ws.Range("I138").Value = 1212.2903
$ws.Range("J138").Value = 2466.5652
$ws.Range("K138").Value = 3636.8709
$ws.Range("L138").Value = 7399.6956
$ws.Range("M138").Value = 1503.1291
$ws.Range("N138").Value = -17679.6956

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 61
$ws.Range("H61").Value = 2109.3555
$ws.Range("I61").Value = 1962.1282
$ws.Range("K61").Value = 1962.1282
$ws.Range("M61").Value = -1750.1282
# Row 74
$ws.Range("H74").Value = 36685
$ws.Range("I74").Value = 55102.05
$ws.Range("J74").Value = 1692.6
$ws.Range("K74").Value = 55102.05
$ws.Range("L74").Value = 1692.6
$ws.Range("M74").Value = -54228.05
$ws.Range("N74").Value = -3440.6
# Row 77
$ws.Range("H77").Value = 36685
$ws.Range("I77").Value = 55102.05
$ws.Range("J77").Value = 1692.6
$ws.Range("K77").Value = 275510.25
$ws.Range("L77").Value = 8463
$ws.Range("M77").Value = -271142.25
$ws.Range("N77").Value = -17199
# Row 102
$ws.Range("H102").Value = 1425
$ws.Range("I102").Value = 1425
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1425
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 197
$ws.Range("N102").ClearContents()
# Row 122
$ws.Range("H122").Value = 1000.96
$ws.Range("I122").Value = 957.8823
$ws.Range("J122").Value = 1092.5
$ws.Range("K122").Value = 2873.6469
$ws.Range("L122").Value = 3277.5
$ws.Range("M122").Value = -423.6468999999997
$ws.Range("N122").Value = -8177.5
# Row 132
$ws.Range("H132").Value = 235571.69
$ws.Range("I132").Value = 41254.566
$ws.Range("J132").Value = 559433.5600000001
$ws.Range("K132").Value = 123763.698
$ws.Range("L132").Value = 1678300.68
$ws.Range("M132").Value = -121233.698
$ws.Range("N132").Value = -1683360.68
# Row 136
$ws.Range("H136").Value = 2109.3555
$ws.Range("I136").Value = 1962.1282
$ws.Range("K136").Value = 5886.3846
$ws.Range("M136").Value = -3336.3846

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 20
$ws.Range("H20").Value = 2230.45
$ws.Range("I20").Value = 2288.889
$ws.Range("K20").Value = 2288.889
$ws.Range("M20").Value = -2041.889
# Row 99
$ws.Range("H99").Value = 1395.4
$ws.Range("I99").Value = 1354.6666
$ws.Range("J99").Value = 1517.6
$ws.Range("K99").Value = 1354.6666
$ws.Range("L99").Value = 1517.6
$ws.Range("M99").Value = 143.3334
$ws.Range("N99").Value = -4513.6

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 25
$ws.Range("H25").Value = 93346.336
$ws.Range("J25").Value = 93346.336
$ws.Range("L25").Value = 93346.336
$ws.Range("N25").Value = -93694.336
# Row 82
$ws.Range("H82").Value = 16500
$ws.Range("J82").Value = 16500
$ws.Range("L82").Value = 16500
$ws.Range("N82").Value = -17222
# Row 85
$ws.Range("H85").Value = 16500
$ws.Range("J85").Value = 16500
$ws.Range("L85").Value = 16500
$ws.Range("N85").Value = -18996
# Row 132
$ws.Range("H132").Value = 49486.24
$ws.Range("I132").Value = 68214.2
$ws.Range("K132").Value = 204642.6
$ws.Range("M132").Value = -202112.6

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 112
$ws.Range("H112").Value = 3305.4
$ws.Range("I112").Value = 1027
$ws.Range("J112").Value = 3875
$ws.Range("K112").Value = 3081
$ws.Range("L112").Value = 11625
$ws.Range("M112").Value = -1973
$ws.Range("N112").Value = -13841
# Row 122
$ws.Range("H122").Value = 730.7037
$ws.Range("I122").Value = 616.8
$ws.Range("J122").Value = 797.7059
$ws.Range("K122").Value = 5551.2
$ws.Range("L122").Value = 7179.3531
$ws.Range("M122").Value = -3101.2
$ws.Range("N122").Value = -12079.3531

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 70
$ws.Range("H70").Value = 4124.5
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 4249
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 4249
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -4789
# Row 73
$ws.Range("H73").Value = 4124.5
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 4249
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 4249
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -6121
# Row 80
$ws.Range("H80").Value = 9193.200000000001
$ws.Range("I80").Value = 2933.111
$ws.Range("J80").Value = 18583.334
$ws.Range("K80").Value = 2933.111
$ws.Range("L80").Value = 18583.334
$ws.Range("M80").Value = -1935.111
$ws.Range("N80").Value = -20579.334
# Row 83
$ws.Range("H83").Value = 9193.200000000001
$ws.Range("I83").Value = 2933.111
$ws.Range("J83").Value = 18583.334
$ws.Range("K83").Value = 14665.555
$ws.Range("L83").Value = 92916.67
$ws.Range("M83").Value = -9673.555
$ws.Range("N83").Value = -102900.67
# Row 97
$ws.Range("H97").Value = 1152.9048
$ws.Range("I97").Value = 1130.7693
$ws.Range("J97").Value = 1188.875
$ws.Range("K97").Value = 1130.7693
$ws.Range("L97").Value = 1188.875
$ws.Range("M97").Value = -634.7692999999999
$ws.Range("N97").Value = -2180.875
# Row 122
$ws.Range("H122").Value = 2081.4583
$ws.Range("I122").Value = 1957.3846
$ws.Range("K122").Value = 5872.1538
$ws.Range("M122").Value = -3422.1538

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 132
$ws.Range("H132").Value = 187765.25
$ws.Range("I132").Value = 42922.36
$ws.Range("K132").Value = 128767.08
$ws.Range("M132").Value = -126237.08
# Row 136
$ws.Range("H136").Value = 501403.4
$ws.Range("I136").Value = 770102.9399999999
$ws.Range("J136").Value = 2390
$ws.Range("K136").Value = 2310308.82
$ws.Range("L136").Value = 7170
$ws.Range("M136").Value = -2307758.82
$ws.Range("N136").Value = -12270

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 96
$ws.Range("H96").Value = 3109.4
$ws.Range("I96").Value = 2500
$ws.Range("J96").Value = 3515.6667
$ws.Range("K96").Value = 2500
$ws.Range("L96").Value = 3515.6667
$ws.Range("M96").Value = -1127
$ws.Range("N96").Value = -6261.6667
# Row 136
$ws.Range("H136").Value = 2037649.6
$ws.Range("I136").Value = 2748704
$ws.Range("J136").Value = 717120
$ws.Range("K136").Value = 8246112
$ws.Range("L136").Value = 2151360
$ws.Range("M136").Value = -8243562
$ws.Range("N136").Value = -2156460
